# Datorama_Trait_Hierarchcial.xlsx — "Added the display level changes to the framework"
#
# Renames the per-sheet "Trait ..." display/label column (column B, rows 3-5 on the
# four TraitDelivery_* sheets, and rows 3-4 on TraitCon_TraitID) to the new
# "... (Segment)" / de-prefixed wording used by the framework going forward:
#
#   Trait Impressions             -> Impressions (Segment)
#   Trait Cost                    -> Media Cost (Segment)
#   Trait Clicks                  -> Clicks (Segment)
#   Trait Click Based Conversions -> Click Based Conversions
#   Trait View Based Conversions  -> View Based Conversions

$wb = $excel.ActiveWorkbook

$newImpressions = "Impressions (Segment)"
$newCost        = "Media Cost (Segment)"
$newClicks      = "Clicks (Segment)"
$newClickConv   = "Click Based Conversions"
$newViewConv    = "View Based Conversions"

# --- TraitDelivery_CampaignID ------------------------------------------------
$ws = $wb.Worksheets.Item("TraitDelivery_CampaignID")
$ws.Range("B3").Value = $newImpressions
$ws.Range("B4").Value = $newCost
$ws.Range("B5").Value = $newClicks
$ws.Range("C5").Select()

# --- TraitDelivery_CampaignTargetID ------------------------------------------
$ws = $wb.Worksheets.Item("TraitDelivery_CampaignTargetID")
$ws.Range("B3").Value = $newImpressions
$ws.Range("B4").Value = $newCost
$ws.Range("B5").Value = $newClicks
$ws.Range("B3:B5").Select()

# --- TraitDelivery_TraitID ---------------------------------------------------
$ws = $wb.Worksheets.Item("TraitDelivery_TraitID")
$ws.Range("B3").Value = $newImpressions
$ws.Range("B4").Value = $newCost
$ws.Range("B5").Value = $newClicks
$ws.Range("B3:B5").Select()

# --- TraitCon_TraitID ---------------------------------------------------------
$ws = $wb.Worksheets.Item("TraitCon_TraitID")
$ws.Range("B3").Value = $newClickConv
$ws.Range("B4").Value = $newViewConv
$ws.Range("C4").Select()

# --- TraitDelivery_AdvertiserID ----------------------------------------------
# Selected/activated last so it ends up as the active tab, matching the saved file.
$ws = $wb.Worksheets.Item("TraitDelivery_AdvertiserID")
$ws.Range("B3").Value = $newImpressions
$ws.Range("B4").Value = $newCost
$ws.Range("B5").Value = $newClicks
$ws.Activate()
$ws.Range("C5").Select()
